# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
# Commit: Updated cryptos list on Fri Nov 15 20:48:09 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) contain numeric-looking text that Excel
# would otherwise auto-convert to numbers; force the cell format to Text first
# so values are written back as literal strings, matching the source XML (t="inlineStr").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.066.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.077.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.57'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.372'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.881'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +13.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.076.19'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.677'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +20.84%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.53%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.802.59'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.36'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.76'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.628.98'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.159.85'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.45'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000221'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.69'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '430.78'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.42'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.90%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.165'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.81'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '512.71'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.84'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.95'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.82%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.136'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.30'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.24%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0721'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.42%  '
$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000275'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +16.95%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '144.51'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.72'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +7.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '165.09'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.16%  '
